$d = $word.ActiveDocument
$pairs = @(
    ,@("[[PERSON_9]] – „od [[PERSON_9]]“", "[[PERSON_9]] – „od [[PERSON_10]]“")
    ,@("[[PERSON_10]] – „pro [[PERSON_10]]“", "[[PERSON_11]] – „pro [[PERSON_11]]“")
    ,@("[[PERSON_11]] – „k [[PERSON_12]]“", "[[PERSON_12]] – „k [[PERSON_12]]“")
    ,@("[[PERSON_23]] – „pro [[PERSON_24]]“", "[[PERSON_23]] – „pro [[PERSON_23]]“")
    ,@("[[PERSON_25]] – „s [[PERSON_25]]“", "[[PERSON_24]] – „s [[PERSON_24]]“")
    ,@("[[PERSON_26]] – „k [[PERSON_26]]“", "[[PERSON_25]] – „k [[PERSON_25]]“")
    ,@("[[PERSON_27]] – „s [[PERSON_27]]“", "[[PERSON_26]] – „s [[PERSON_26]]“")
    ,@("[[PERSON_28]] – „o [[PERSON_29]]“", "[[PERSON_27]] – „o [[PERSON_28]]“")
    ,@("[[PERSON_30]] – „pro [[PERSON_30]]“", "[[PERSON_29]] – „pro [[PERSON_29]]“")
    ,@("[[PERSON_31]] – „s [[PERSON_31]]“", "[[PERSON_30]] – „s [[PERSON_30]]“")
    ,@("[[PERSON_32]] – „k [[PERSON_32]]“", "[[PERSON_31]] – „k [[PERSON_31]]“")
    ,@("[[PERSON_33]] – „s [[PERSON_33]]“", "[[PERSON_32]] – „s [[PERSON_32]]“")
    ,@("[[PERSON_34]] – „o [[PERSON_34]]“", "[[PERSON_33]] – „o [[PERSON_33]]“")
    ,@("[[PERSON_35]] – „u [[PERSON_35]]“", "[[PERSON_34]] – „u [[PERSON_35]]“")
    ,@("[[PERSON_47]] – „o [[PERSON_48]]“", "[[PERSON_47]] – „o [[PERSON_47]]“")
    ,@("[[PERSON_49]] – „k [[PERSON_49]]“", "[[PERSON_48]] – „k [[PERSON_48]]“")
    ,@("V těchto řízeních bylo jednáno např. s [[PERSON_2]], [[PERSON_6]], [[PERSON_28]] či [[PERSON_50]].", "V těchto řízeních bylo jednáno např. s [[PERSON_2]], [[PERSON_6]], [[PERSON_27]] či [[PERSON_49]].")
    ,@("svědek [[PERSON_33]] (ve výpovědi označen jako „svědek Černého“),", "svědek [[PERSON_32]] (ve výpovědi označen jako „svědek Černého“),")
    ,@("tlumočník [[PERSON_35]], zapsaný v seznamu tlumočníků.", "tlumočník [[PERSON_34]], zapsaný v seznamu tlumočníků.")
    ,@("Oční vyšetření č. OFT/2023/11281 provedené MUDr. [[PERSON_34]].", "Oční vyšetření č. OFT/2023/11281 provedené MUDr. [[PERSON_33]].")
    ,@("právní cloud účet ID: LEX-ACC-88221 (spravovala [[PERSON_49]]),", "právní cloud účet ID: LEX-ACC-88221 (spravovala [[PERSON_48]]),")
    ,@("PhDr. [[PERSON_34]] – psychologický posudek,", "PhDr. [[PERSON_33]] – psychologický posudek,")
    ,@("MUDr. [[PERSON_28]] – posudek z traumatologie,", "MUDr. [[PERSON_27]] – posudek z traumatologie,")
    ,@("Ing. [[PERSON_9]] – expertiza IT infrastruktury.", "Ing. [[PERSON_10]] – expertiza IT infrastruktury.")
    ,@("Tyto účty byly doloženy např. od [[PERSON_23]], [[PERSON_42]] nebo [[PERSON_31]].", "Tyto účty byly doloženy např. od [[PERSON_23]], [[PERSON_42]] nebo [[PERSON_30]].")
    ,@("[[PERSON_50]],", "[[PERSON_49]],")
    ,@("[[PERSON_27]],", "[[PERSON_26]],")
    ,@("[[PERSON_11]].", "[[PERSON_12]].")
)

$failures = 0
foreach ($pair in $pairs) {
    $oldText = $pair[0]
    $newText = $pair[1]
    $rng = $d.Content
    $found = $rng.Find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        $failures++
        Write-Host "NOT FOUND: $oldText"
    }
}
Write-Host "Done. Failures: $failures"
